$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: add Field Values (column E), wrap text, and increase row height
$ws.Range("E8").Value = "username = COSME0007`n  "
$ws.Range("E8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 31.5

# Row 14
$ws.Range("C14").Value = "SME Customer Prospect"
$ws.Range("D14").Value = "SME Costomer Prospects"
$ws.Range("E14").Value = "prospectId = COSME0013"
$ws.Range("F14").Value = "COSME0013COSMECUSTOMERRequest"
$ws.Range("G14").Value = "COSME0013SearchSMEProspectResponse"

# Row 16
$ws.Range("C16").Value = "SMEdatalist"
$ws.Range("D16").Value = "SME Datalist"
$ws.Range("E16").Value = "code = 1234"
$ws.Range("F16").Value = "COSME0015GetSMERequest"
$ws.Range("G16").Value = "COSME0015GetSMEResponse"

# Row 17
$ws.Range("C17").Value = "OnBoardSMECustomer"
$ws.Range("D17").Value = "OnBoardSMECustomer"
$ws.Range("E17").Value = "prospectId = COSME0016"
$ws.Range("F17").Value = "COSME0016PostonboardSMEcustomerRequest"
$ws.Range("G17").Value = "COSME0016PostonboardSMEcustomerResponse"

# Row 18
$ws.Range("C18").Value = "DocumentsUploadRequest"
$ws.Range("D18").Value = "DocumentsUploadRequest"
$ws.Range("E18").Value = "fileName = COSME0017"
$ws.Range("F18").Value = "COSME0017DocumentsUploadRequest"
$ws.Range("G18").Value = "COSME0017DocumentsUploadResponse"

# Row 19
$ws.Range("C19").Value = "GETDocumentBy ID"
$ws.Range("D19").Value = "GETDocumentBy ID"
$ws.Range("E19").Value = "prospectId =  COSME0018"
$ws.Range("F19").Value = "COSME0018GetDocumentbyIdRequest"
$ws.Range("G19").Value = "COSME0018GetDocumentbyidResponse"

# Row 20
$ws.Range("C20").Value = "ReuploadDocument"
$ws.Range("D20").Value = "ReuploadDocument"
$ws.Range("E20").Value = "prospectId = COSME0019"
$ws.Range("F20").Value = "COMSE0019PUTReuploadDocumentRequest"
$ws.Range("G20").Value = "COSME0017DocumentsUploadResponse"

# Row 21
$ws.Range("C21").Value = "PreScreening "
$ws.Range("D21").Value = "PreScreening "
$ws.Range("E21").Value = "prospectId = COSME0020"
$ws.Range("F21").Value = "COSME0020GETPreScreenRequest"
$ws.Range("G21").Value = "COSME0020GETPreScreenResponse"

# sheetView changes: topLeftCell and selection
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F14").Select()
